$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Configure the new token's modifiers (player or npc row)
$ws.Range("B3").Value = "Hand Seal"
$ws.Range("C3").Value = "All"
$ws.Range("D3").Value = "All"
$ws.Range("E3").Value = "Roll Add"
$ws.Range("F3").Value = "All"

# Move the active selection to the computed code cell
$ws.Activate()
$ws.Range("I4").Select()
